$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old sub-header row (Hiver / Eté / Année / unit labels) - this
# shifts all data rows up by one and updates the dimension automatically.
$ws.Rows.Item(2).Delete()

# Rewrite row 1 as the new header row: idx, idx2, Name, Date Start, Date End,
# (m3/s), (MW1), (MW2), (GWh) Winter, (GWh) Summer, (GWh) Year
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# Give the unit-style headers (F1:K1) the same look as the rest of row 1
# (Arial 9pt, non-bold/non-italic/non-underlined) without imposing any
# particular number format on them.
$ws.Range("F1:K1").Font.Name = "Arial"
$ws.Range("F1:K1").Font.Size = 9
$ws.Range("F1:K1").Font.Bold = $false
$ws.Range("F1:K1").Font.Italic = $false
$ws.Range("F1:K1").Font.Underline = $false

# Restore the default selection seen in the final workbook.
$ws.Range("G19").Select()
